$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the watering/fill rate values for each flower type (column C)
$ws.Range("C5").Value = 50
$ws.Range("C8").Value = 45
$ws.Range("C11").Value = 40
$ws.Range("C14").Value = 35
$ws.Range("C17").Value = 35
$ws.Range("C20").Value = 30

# Update the selected cell in the sheet view
$ws.Range("D22").Select()
